$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the comma-separated "Recorded By" names in column G for the affected rows.
# Same underlying sets of people per row; only the ordering changes (per diff).
$ws.Cells.Item(2, 7).Value = "Administrator, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Veronia Rafat, Dr. Amira Sobhy"
$ws.Cells.Item(3, 7).Value = "Dr. Majorelle Magdy, Dr. Eman Tantawi, Administrator, Dr. Asmaa Reda, Dr. Veronia Rafat, Dr. Hend Mahmoud"
$ws.Cells.Item(4, 7).Value = "Dr. Majorelle Magdy, Dr. Eman Tantawi, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Asmaa Reda, Dr. Hend Mahmoud"
$ws.Cells.Item(5, 7).Value = "Dr. Eman Tantawi, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Asmaa Reda"
$ws.Cells.Item(6, 7).Value = "Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Menna tuâ€™Allah Medhat, Dr. Manar Montaser, Dr. Alshimaa Atef"
$ws.Cells.Item(7, 7).Value = "Dr. Amera Ahmad Saad, Dr. Abeer Ragab, Dr. Lamiaa Ossama, Dr. Menna tu'Alllah Mohammad, Dr. Nada Mohammad, Dr. Kerelos Zareef, Dr. Fatma Elhady"
$ws.Cells.Item(8, 7).Value = "Dr. Nada Mohammad, Dr. Abeer Ragab"
$ws.Cells.Item(11, 7).Value = "Dr. Amal Awwad, Dr. Aya Saeed, Dr. Safa Hany"
$ws.Cells.Item(12, 7).Value = "Dr. Amira Ibrahim, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Marina Youhanna"
$ws.Cells.Item(13, 7).Value = "Dr. Esraa Mostafa, Dr. Amira Ibrahim, Dr. Yasmeena Fattoh"
$ws.Cells.Item(15, 7).Value = "Dr. Mohammad Safwat, Dr. Rania Ahmad Youssef"
$ws.Cells.Item(20, 7).Value = "Dr. Mohammad Safwat, Dr. Mariam Toma Gerges"
$ws.Cells.Item(25, 7).Value = "Menna tuâ€™Allah Gamil, Dr. Nouran Mahmoud"
$ws.Cells.Item(27, 7).Value = "Dr. Hana Amr, Dr. Nourham Mostafa"
$ws.Cells.Item(30, 7).Value = "Dr. Aya Hanafy, Dr. Yassmen Ahmad, Dr. Shorok Mohammad, Dr. Wafaa Ebida"
